# Weekly update: insert a new price observation as the new first data
# row for this subset (row 391), pushing all existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 391; this shifts rows 391:456 down to 392:457
# and extends the sheet dimension to A1:R457.
$ws.Rows.Item(391).Insert()

# Populate the newly inserted row with this week's data.
$ws.Cells.Item(391, 1).Value  = 8
$ws.Cells.Item(391, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(391, 3).Value  = "Coquimbo"
$ws.Cells.Item(391, 4).Value  = 44951
$ws.Cells.Item(391, 5).Value  = 4
$ws.Cells.Item(391, 6).Value  = 100112032
$ws.Cells.Item(391, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(391, 8).Value  = "Sin especificar"
$ws.Cells.Item(391, 9).Value  = "Primera"
$ws.Cells.Item(391, 10).Value = 400
$ws.Cells.Item(391, 11).Value = 9000
$ws.Cells.Item(391, 12).Value = 10000
$ws.Cells.Item(391, 13).Value = 9500
$ws.Cells.Item(391, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(391, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(391, 16).Value = 136
$ws.Cells.Item(391, 17).Value = 70
$ws.Cells.Item(391, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date/time display format used
# throughout column D.
$ws.Cells.Item(391, 4).NumberFormat = $ws.Cells.Item(392, 4).NumberFormat
